$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-20 Wednesday" "2024-03-21 Thursday"

Replace-Text "752÷8=94, 0" "411÷7=58, 5"
Replace-Text "698÷4=174, 2" "341÷6=56, 5"
Replace-Text "298÷3=99, 1" "178÷8=22, 2"
Replace-Text "719÷6=119, 5" "697÷8=87, 1"
Replace-Text "928÷8=116, 0" "642÷4=160, 2"

Replace-Text "509÷4=127, 1" "878÷8=109, 6"
Replace-Text "785÷6=130, 5" "766÷2=383, 0"
Replace-Text "356÷3=118, 2" "841÷7=120, 1"
Replace-Text "164÷7=23, 3" "795÷5=159, 0"
Replace-Text "932÷7=133, 1" "916÷5=183, 1"

Replace-Text "282÷9=31, 3" "964÷2=482, 0"
Replace-Text "120÷6=20, 0" "932÷6=155, 2"
Replace-Text "887÷4=221, 3" "672÷5=134, 2"
Replace-Text "125÷3=41, 2" "496÷6=82, 4"
Replace-Text "557÷9=61, 8" "738÷7=105, 3"

Replace-Text "419÷4=104, 3" "210÷3=70, 0"
Replace-Text "750÷5=150, 0" "751÷3=250, 1"
Replace-Text "203÷4=50, 3" "305÷7=43, 4"
Replace-Text "184÷3=61, 1" "495÷6=82, 3"
Replace-Text "176÷7=25, 1" "338÷3=112, 2"

Replace-Text "871÷5=174, 1" "948÷2=474, 0"
Replace-Text "667÷5=133, 2" "750÷7=107, 1"
Replace-Text "536÷5=107, 1" "322÷2=161, 0"
Replace-Text "654÷8=81, 6" "923÷6=153, 5"
Replace-Text "744÷6=124, 0" "176÷9=19, 5"
